# Commit: "unify the conception of DataNode, DataTable, Entity."
# The sheet formerly called "Property1" is renamed to "DataNode" to match
# the new unified naming convention, and the in-sheet selection is left
# parked at C36 (matching the state the workbook was saved in).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: "Property1" -> "DataNode"
$ws.Name = "DataNode"

# Move / leave the active selection at C36, as captured in the saved file.
$ws.Range("C36").Select()
